# Ready for production (just for now)
#
# Updates the "Plan1" course-code table:
#  - rows 57-60 get new sequential course codes (106-109) and their
#    descriptions are replaced with combined-course labels
#  - seven new rows (61-67) are appended: three more combined-course rows,
#    then the original PGMAT/PROFMAT/PROFMAT(Verão)/Honors rows (which used
#    to live at 57-60) are moved down to 64-67, keeping their old codes
#  - the sheet selection is moved to the newly added tail of the table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 57; Code = 106; Text = "Engenharia Civil, Ciência da Computação" },
    @{ Row = 58; Code = 107; Text = "Licenciatura em Química Noturno, Matemática Industrial" },
    @{ Row = 59; Code = 108; Text = "Física Noturno, Matemática Industrial" },
    @{ Row = 60; Code = 109; Text = "Matemática Diurno, Matemática Industrial" },
    @{ Row = 61; Code = 110; Text = "Matemática Diurno, Matemática Industrial" },
    @{ Row = 62; Code = 111; Text = "Matemática Diurno, Matemática Industrial" },
    @{ Row = 63; Code = 112; Text = "Matemática Industrial, Engenharia Industrial Madeireira Diurno" },
    @{ Row = 64; Code = 42;  Text = "PGMAT" },
    @{ Row = 65; Code = 43;  Text = "PROFMAT" },
    @{ Row = 66; Code = 44;  Text = "PROFMAT(Verão)" },
    @{ Row = 67; Code = 46;  Text = "Honors" }
)

foreach ($u in $updates) {
    $a = $ws.Cells.Item($u.Row, 1)
    $a.Value = $u.Code
    $a.NumberFormat = "000"
    $a.Font.Bold = $true

    $b = $ws.Cells.Item($u.Row, 2)
    $b.Value = $u.Text
}

$ws.Range("A64").Select()

Write-Output "done"
